$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '67.751.55'
$ws.Range('E2').Value = '  +1.50%  '

# Row 3
$ws.Range('D3').Value = '3.124.19'
$ws.Range('E3').Value = '  +1.88%  '

# Row 4
$ws.Range('E4').Value = '  -0.10%  '

# Row 5
$ws.Range('D5').Value = '''577.10'
$ws.Range('E5').Value = '  +0.19%  '

# Row 6
$ws.Range('D6').Value = '''179.94'
$ws.Range('E6').Value = '  +7.26%  '

# Row 7
$ws.Range('D7').Value = '''0.999'
$ws.Range('E7').Value = '  -0.10%  '

# Row 8
$ws.Range('D8').Value = '3.123.60'
$ws.Range('E8').Value = '  +1.97%  '

# Row 9
$ws.Range('D9').Value = '''0.518'
$ws.Range('E9').Value = '  +1.19%  '

# Row 10
$ws.Range('D10').Value = '''6.53'
$ws.Range('E10').Value = '  +2.32%  '

# Row 11
$ws.Range('E11').Value = '  +2.00%  '

# Row 12
$ws.Range('E12').Value = '  +0.06%  '

# Row 13
$ws.Range('D13').Value = '''0.0000242'
$ws.Range('E13').Value = '  +1.02%  '

# Row 14
$ws.Range('D14').Value = '''36.81'
$ws.Range('E14').Value = '  +3.09%  '

# Row 15
$ws.Range('E15').Value = '  +0.83%  '

# Row 16
$ws.Range('D16').Value = '3.641.33'
$ws.Range('E16').Value = '  +1.72%  '

# Row 17
$ws.Range('D17').Value = '67.618.97'
$ws.Range('E17').Value = '  +1.28%  '

# Row 18
$ws.Range('D18').Value = '''7.05'
$ws.Range('E18').Value = '  +0.94%  '

# Row 19
$ws.Range('D19').Value = '3.118.60'
$ws.Range('E19').Value = '  +1.51%  '

# Row 20
$ws.Range('D20').Value = '''16.45'
$ws.Range('E20').Value = '  -2.02%  '

# Row 21
$ws.Range('D21').Value = '''487.63'
$ws.Range('E21').Value = '  -0.15%  '

# Row 22
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').Value = '''7.77'
$ws.Range('E22').Value = '  +1.11%  '

# Row 23
$ws.Range('B23').Value = 'Polygon'
$ws.Range('C23').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D23').Value = '''0.692'
$ws.Range('E23').Value = '  +0.58%  '

# Row 24
$ws.Range('D24').Value = '''83.89'
$ws.Range('E24').Value = '  +1.38%  '

# Row 25
$ws.Range('B25').Value = 'Fetch.AI'
$ws.Range('C25').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D25').Value = '''2.32'
$ws.Range('E25').Value = '  +5.34%  '

# Row 26
$ws.Range('B26').Value = 'InternetComputer(DFINITY)'
$ws.Range('C26').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D26').Value = '''12.82'
$ws.Range('E26').Value = '  +1.19%  '

# Row 27
$ws.Range('D27').Value = '''10.66'
$ws.Range('E27').Value = '  +4.87%  '

# Row 28
$ws.Range('E28').Value = '  +0.04%  '

# Row 29
$ws.Range('D29').Value = '''8.12'
$ws.Range('E29').Value = '  +4.60%  '

# Row 30
$ws.Range('D30').Value = '''2.36'
$ws.Range('E30').Value = '  +4.23%  '

# Row 31
$ws.Range('D31').Value = '''2.63'
$ws.Range('E31').Value = '  +0.64%  '

# Row 32
$ws.Range('D32').Value = '''28.20'
$ws.Range('E32').Value = '  +2.33%  '

# Row 33
$ws.Range('E33').Value = '  +1.35%  '

# Row 34
$ws.Range('D34').Value = '0.0₃0950'
$ws.Range('E34').Value = '  +4.16%  '

# Row 35
$ws.Range('D35').Value = '''0.999'
$ws.Range('E35').Value = '  -0.08%  '

# Row 36
$ws.Range('D36').Value = '''48.63'
$ws.Range('E36').Value = '  +4.13%  '

# Row 37
$ws.Range('D37').Value = '''0.952'
$ws.Range('E37').Value = '  +0.43%  '

# Row 38
$ws.Range('E38').Value = '  -0.18%  '

# Row 39
$ws.Range('E39').Value = '  +7.45%  '

# Row 40
$ws.Range('E40').Value = '  +3.20%  '

# Row 41
$ws.Range('D41').Value = '''49.25'
$ws.Range('E41').Value = '  +0.26%  '

# Row 42
$ws.Range('E42').Value = '  +1.25%  '

# Row 43
$ws.Range('D43').Value = '''8.35'
$ws.Range('E43').Value = '  +0.46%  '

# Row 44
$ws.Range('D44').Value = '''2.69'
$ws.Range('E44').Value = '  +8.67%  '

# Row 45
$ws.Range('D45').Value = '2.785.98'
$ws.Range('E45').Value = '  +1.15%  '

# Row 46
$ws.Range('D46').Value = '''379.71'
$ws.Range('E46').Value = '  +2.88%  '

# Row 47
$ws.Range('E47').Value = '  +0.72%  '

# Row 48
$ws.Range('D48').Value = '''26.79'
$ws.Range('E48').Value = '  +10.07%  '

# Row 49
$ws.Range('D49').Value = '''135.97'
$ws.Range('E49').Value = '  +0.00%  '

# Row 50
$ws.Range('E50').Value = '  +0.02%  '

# Row 51
$ws.Range('D51').Value = '''2.34'
$ws.Range('E51').Value = '  +8.55%  '
